$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is removed entirely; all subsequent
# columns (reviews_average, latitude, longitude, is_permanently_closed,
# gmaps_link, latest_review_date) shift one column to the left.
$ws.Range("E:E").EntireColumn.Delete()
